$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (AD:AF) by cloning the formatting of the
# last existing header cell (AC1, style "1": bold + border + center/top
# alignment) into the new header cells, then overwrite their text.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
# Every row in this sheet shares the same season record: 90-72-1.
for ($row = 2; $row -le 41; $row++) {
    $ws.Cells.Item($row, 30).Value = 90
    $ws.Cells.Item($row, 31).Value = 72
    $ws.Cells.Item($row, 32).Value = 1
}

Write-Output "done"
